$p = $ppt.ActivePresentation

# --- Slide 2 ("Final Phase objectives" title) ---------------------------
# The title text used to be split into three runs: "Final", " " and
# "Phase objectives". The first two runs get merged into a single
# "Final " run (the trailing space moves into the first run).
$s2 = $p.Slides.Item(2)
$titleShape = $s2.Shapes.Item("Titel 1")
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Characters(1, 6).Text = "Final "

# --- Slide 5 ("Experiments" body) ---------------------------------------
# "Made experiments by modifying parameters, then measured the " + "result"
# becomes three runs: "...measured ", "the " and "result." (period added).
$s5 = $p.Slides.Item(5)
$bodyShape = $s5.Shapes.Item("Tijdelijke aanduiding voor inhoud 3")
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Characters(57, 4).Text = "the "
$bodyRange.Characters(61, 7).Text = "result."
